$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 44

$ws.Cells.Item($row, 1).Value = "2025-08-22 09:39:48 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-22 15:09:48 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = ""

# Match formatting used by the rest of the log rows (centered alignment)
$rowRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 8))
$rowRange.HorizontalAlignment = -4108
$rowRange.VerticalAlignment = -4108
